# Update "produtos" sheet: add Qnt-Normal/Qnt-FULL columns (now E/F), keep
# Media-Preco (G), add Mediana-Preco (H), add Media-Vendas (I) and
# Mediana-Vendas (J). Refresh scrapy_datetime + all computed metrics for
# every product row with the freshly scraped results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----------------------------------------------
$ws.Cells.Item(1, 5).Value = "Qnt-Normal"
$ws.Cells.Item(1, 6).Value = "Qnt-FULL"
$ws.Cells.Item(1, 7).Value = "Media-Preco"
$ws.Cells.Item(1, 8).Value = "Mediana-Preco"

# New header cells I1/J1: copy the formatting of an existing bold/bordered
# header cell so the new columns look the same as the rest of row 1.
$ws.Cells.Item(1, 4).Copy() | Out-Null
$ws.Cells.Item(1, 9).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 10).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(1, 9).Value = "Media-Vendas"
$ws.Cells.Item(1, 10).Value = "Mediana-Vendas"

# ---- Data rows (2-11) --------------------------------------------------
$datetime = "2022-05-20 16:23:08"

$data = @(
  @{ Row = 2;  E = "17.105 resultado";  F = "233 resultado";   G = 93;                H = 78;   I = 1491.157894736842; J = 1235 },
  @{ Row = 3;  E = "10.821 resultado";  F = "44 resultado";    G = 126.8235294117647;  H = 135;  I = 89.76470588235294; J = 56 },
  @{ Row = 4;  E = "600 resultado";     F = "74 resultado";    G = 66.31578947368421;  H = 44;   I = 585.421052631579;  J = 460 },
  @{ Row = 5;  E = "2.788 resultado";   F = "42 resultado";    G = 253.3571428571429;  H = 172;  I = 310.4285714285714; J = 227.5 },
  @{ Row = 6;  E = "18.756 resultado";  F = "327 resultado";   G = 55.63157894736842;  H = 40;   I = 1968.526315789474; J = 1337 },
  @{ Row = 7;  E = "1.108 resultado";   F = "52 resultado";    G = 196.7222222222222;  H = 139;  I = 29;                J = 9.5 },
  @{ Row = 8;  E = "2.382 resultado";   F = "36 resultado";    G = 64.07692307692308;  H = 54;   I = 1541.153846153846; J = 1108 },
  @{ Row = 9;  E = "150.268 resultado"; F = "3.166 resultado"; G = 61.68421052631579;  H = 48;   I = 14095.57894736842; J = 6781 },
  @{ Row = 10; E = "679 resultado";     F = "NaoTem";          G = 588.1666666666666;  H = 569;  I = 180.5;             J = 97.5 },
  @{ Row = 11; E = "38.178 resultado";  F = "500 resultado";   G = 78.16666666666667;  H = 47.5; I = 2699.666666666667; J = 2332.5 }
)

foreach ($item in $data) {
  $r = $item.Row
  $ws.Cells.Item($r, 4).Value = $datetime
  $ws.Cells.Item($r, 5).Value = $item.E
  $ws.Cells.Item($r, 6).Value = $item.F
  $ws.Cells.Item($r, 7).Value = $item.G
  $ws.Cells.Item($r, 8).Value = $item.H
  $ws.Cells.Item($r, 9).Value = $item.I
  $ws.Cells.Item($r, 10).Value = $item.J
}
